$wb = $excel.ActiveWorkbook

# --- NewLoanInput sheet -----------------------------------------------------
# Previously the active/selected sheet (tabSelected) with B2 selected.
# The edit moves the selection to B5 and the sheet is no longer the active
# tab (we activate "Repayment schedule" afterwards, which becomes the new
# active tab).
$wsInput = $wb.Worksheets.Item("NewLoanInput")
$wsInput.Activate() | Out-Null
$wsInput.Range("B5").Select() | Out-Null

# --- Repayment schedule sheet ------------------------------------------------
# Insert a new (blank) column at N, shifting the old N:P ("Late",
# "heading"/Original, "Outstanding") data right to O:Q. This mirrors a
# "Variable Instalments" schedule column being added ahead of the existing
# Late/Original/Outstanding columns.
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Columns("N").Insert() | Out-Null

# The newly inserted column inherits the column-M formatting by default in
# Excel; set its width to match the recorded width (stored width 11).
$wsRepay.Columns("N").ColumnWidth = 10.166666666666666

# This sheet becomes the active sheet/tab, with R10 selected.
$wsRepay.Activate() | Out-Null
$wsRepay.Range("R10").Select() | Out-Null
